$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the old header cell value (C6 keeps its style, but text is removed)
$ws.Range("C6").Value = $null
$ws.Range("D6").Value = $null
$ws.Range("E6").Value = $null

# Row 7: spacious / geniş, ferah / sıfat
$ws.Range("C7").Value = "spacious"
$ws.Range("D7").Value = "geniş, ferah"
$ws.Range("E7").Value = "sıfat"

# Row 8: legibility / okunabilirlik / isim
$ws.Range("C8").Value = "legibility"
$ws.Range("E8").Value = "isim"
$ws.Range("D8").Value = "okunabilirlik"

# Update selection to C9
$ws.Range("C9").Select() | Out-Null
